$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old SupId/CategoryName columns (C, D, I) that are no longer used
$ws.Range("C4:D5").Value = $null
$ws.Range("I4:I5").Value = $null

# New header row (row 4), columns E:H
$ws.Range("E4").Value = "SupName"
$ws.Range("F4").Value = "Description"
$ws.Range("G4").Value = "Uom"
$ws.Range("H4").Value = "Price"

# Row 5
$ws.Range("E5").Value = "ALPHA OFFICE Supplies"
$ws.Range("F5").Value = "Exercise Book (10 pg)"
$ws.Range("G5").Value = "Dozen"
$ws.Range("H5").Value = 50

# Write the "Exercise Book (75 pg)" description (row 7) before the row 6
# values so the shared-string table ends up in the same first-use order
# as the original authored file.
$ws.Range("F7").Value = "Exercise Book (75 pg)"

# Row 6
$ws.Range("E6").Value = "ALPHA OFFICE Supplies"
$ws.Range("F6").Value = "Envelope White (5`"x7`")"
$ws.Range("G6").Value = "Each"
$ws.Range("H6").Value = 50

# Row 7 (remaining cells)
$ws.Range("E7").Value = "ALPHA OFFICE Supplies"
$ws.Range("G7").Value = "Dozen"
$ws.Range("H7").Value = 60

# Column widths: closest achievable values to the target widths given engine quantization
$ws.Columns.Item(4).ColumnWidth = 16.666666666666668
$ws.Columns.Item(5).ColumnWidth = 19.833333333333332
$ws.Columns.Item(6).ColumnWidth = 22.333333333333332

# Selection
$ws.Range("H6").Select()
